# Update computed profit figures in the Midgardsormr_Profits crafting-leve sheets
# (scheduled price-refresh run).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1120.8
$ws.Range("I15").Value = 1120.8
$ws.Range("K15").Value = 3362.4
$ws.Range("M15").Value = -3193.4
$ws.Range("H55").Value = 168.33333
$ws.Range("I55").Value = 172.5
$ws.Range("J55").Value = 160
$ws.Range("K55").Value = 172.5
$ws.Range("L55").Value = 160
$ws.Range("M55").Value = 41.5
$ws.Range("N55").Value = -588
$ws.Range("H107").Value = 1078.9166
$ws.Range("I107").Value = 1158.8182
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 1158.8182
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = 761.1818000000001
$ws.Range("N107").Value = -4040
$ws.Range("H111").Value = 2002
$ws.Range("I111").Value = 1915.375
$ws.Range("J111").Value = 2175.25
$ws.Range("K111").Value = 5746.125
$ws.Range("L111").Value = 6525.75
$ws.Range("M111").Value = -2679.125
$ws.Range("N111").Value = -12659.75
$ws.Range("H127").Value = 1304.7894
$ws.Range("J127").Value = 1953.6
$ws.Range("L127").Value = 5860.799999999999
$ws.Range("N127").Value = -15780.8
$ws.Range("H137").Value = 32889.85
$ws.Range("I137").Value = 77560.25
$ws.Range("J137").Value = 21722.25
$ws.Range("K137").Value = 232680.75
$ws.Range("L137").Value = 65166.75
$ws.Range("M137").Value = -230130.75
$ws.Range("N137").Value = -70266.75
$ws.Range("H138").Value = 25951.555
$ws.Range("I138").Value = 3077.889
$ws.Range("J138").Value = 41200.668
$ws.Range("K138").Value = 9233.667000000001
$ws.Range("L138").Value = 123602.004
$ws.Range("M138").Value = -4093.667000000001
$ws.Range("N138").Value = -133882.004

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H32").Value = 16508.184
$ws.Range("I32").Value = 18519.967
$ws.Range("K32").Value = 18519.967
$ws.Range("M32").Value = -18232.967
$ws.Range("H132").Value = 2290.0862
$ws.Range("I132").Value = 2063.625
$ws.Range("J132").Value = 3377.1
$ws.Range("K132").Value = 6190.875
$ws.Range("L132").Value = 10131.3
$ws.Range("M132").Value = -3660.875
$ws.Range("N132").Value = -15191.3

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 5000
$ws.Range("I8").Value = 5000
$ws.Range("K8").Value = 5000
$ws.Range("M8").Value = -4860
$ws.Range("H20").Value = 50376
$ws.Range("I20").Value = 79842.875
$ws.Range("K20").Value = 79842.875
$ws.Range("M20").Value = -79595.875
$ws.Range("H82").Value = 16045.333
$ws.Range("H85").Value = 16045.333
$ws.Range("H86").Value = 1613.05
$ws.Range("I86").Value = 1531.2174
$ws.Range("J86").Value = 1723.7646
$ws.Range("K86").Value = 1531.2174
$ws.Range("L86").Value = 1723.7646
$ws.Range("M86").Value = -408.2174
$ws.Range("N86").Value = -3969.7646
$ws.Range("H89").Value = 1613.05
$ws.Range("I89").Value = 1531.2174
$ws.Range("J89").Value = 1723.7646
$ws.Range("K89").Value = 7656.087
$ws.Range("L89").Value = 8618.823
$ws.Range("M89").Value = -2040.087
$ws.Range("N89").Value = -19850.823
$ws.Range("H94").Value = 3031071
$ws.Range("I94").Value = 747.7083
$ws.Range("K94").Value = 747.7083
$ws.Range("M94").Value = -296.7083
$ws.Range("H99").Value = 1807.3871
$ws.Range("I99").Value = 1722.96
$ws.Range("K99").Value = 1722.96
$ws.Range("M99").Value = -224.96

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 7500
$ws.Range("I10").Value = 5000
$ws.Range("J10").Value = 10000
$ws.Range("K10").Value = 5000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = -4861
$ws.Range("N10").Value = -10278
$ws.Range("H19").Value = 21200850
$ws.Range("I19").Value = 31800624
$ws.Range("J19").Value = 1300
$ws.Range("K19").Value = 31800624
$ws.Range("L19").Value = 1300
$ws.Range("M19").Value = -31800454
$ws.Range("N19").Value = -1640
$ws.Range("H22").Value = 283.26666
$ws.Range("I22").Value = 195.11111
$ws.Range("K22").Value = 195.11111
$ws.Range("M22").Value = 154.88889
$ws.Range("H24").Value = 21200850
$ws.Range("I24").Value = 31800624
$ws.Range("J24").Value = 1300
$ws.Range("K24").Value = 31800624
$ws.Range("L24").Value = 1300
$ws.Range("M24").Value = -31800454
$ws.Range("N24").Value = -1640
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H31").Value = 2328933.8
$ws.Range("I31").Value = 7146421
$ws.Range("J31").Value = 3250.3103
$ws.Range("K31").Value = 7146421
$ws.Range("L31").Value = 3250.3103
$ws.Range("M31").Value = -7146126
$ws.Range("N31").Value = -3840.3103
$ws.Range("H34").Value = 2328933.8
$ws.Range("I34").Value = 7146421
$ws.Range("J34").Value = 3250.3103
$ws.Range("K34").Value = 7146421
$ws.Range("L34").Value = 3250.3103
$ws.Range("M34").Value = -7146219
$ws.Range("N34").Value = -3654.3103
$ws.Range("H86").Value = 39840.42
$ws.Range("I86").Value = 49070.07
$ws.Range("J86").Value = 13997.4
$ws.Range("K86").Value = 49070.07
$ws.Range("L86").Value = 13997.4
$ws.Range("M86").Value = -47947.07
$ws.Range("N86").Value = -16243.4
$ws.Range("H89").Value = 39840.42
$ws.Range("I89").Value = 49070.07
$ws.Range("J89").Value = 13997.4
$ws.Range("K89").Value = 245350.35
$ws.Range("L89").Value = 69987
$ws.Range("M89").Value = -239734.35
$ws.Range("N89").Value = -81219
$ws.Range("H132").Value = 113656.22
$ws.Range("I132").Value = 144700.86
$ws.Range("K132").Value = 434102.58
$ws.Range("M132").Value = -431572.58

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 640
$ws.Range("I134").Value = 640
$ws.Range("K134").Value = 1920
$ws.Range("M134").Value = 3150
$ws.Range("H139").Value = 4923.9585
$ws.Range("I139").Value = 3908.75
$ws.Range("K139").Value = 11726.25
$ws.Range("M139").Value = -6586.25

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3673.5
$ws.Range("I40").Value = 3078.6
$ws.Range("K40").Value = 3078.6
$ws.Range("M40").Value = -2942.6
$ws.Range("H55").Value = 1771.5312
$ws.Range("I55").Value = 1365
$ws.Range("J55").Value = 2130.2354
$ws.Range("K55").Value = 1365
$ws.Range("L55").Value = 2130.2354
$ws.Range("M55").Value = -1192
$ws.Range("N55").Value = -2476.2354
$ws.Range("H100").Value = 1820.9166
$ws.Range("J100").Value = 2326.3333
$ws.Range("L100").Value = 2326.3333
$ws.Range("N100").Value = -3408.3333

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 50000
$ws.Range("J40").Value = 50000
$ws.Range("L40").Value = 50000
$ws.Range("N40").Value = -50298
